$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = "261.64"
    "E2"  = "1.45%"
    "E3"  = "1.33%"
    "D4"  = "4.771"
    "E4"  = "2.16%"
    "D5"  = "0.06072"
    "E5"  = "2.77%"
    "D6"  = "6.721"
    "E6"  = "1.20%"
    "D7"  = "0.8648"
    "E7"  = "1.23%"
    "D8"  = "0.9241"
    "E8"  = "-2.57%"
    "D9"  = "0.1412"
    "E9"  = "0.74%"
    "D10" = "0.05002"
    "E10" = "-0.03%"
    "D11" = "0.07154"
    "E11" = "0.88%"
    "D12" = "0.03063"
    "E12" = "-1.43%"
    "D13" = "0.09114"
    "E13" = "-0.42%"
    "D14" = "0.001529"
    "E14" = "0.35%"
    "D15" = "0.0006101"
    "E15" = "1.18%"
    "D16" = "0.006195"
    "E16" = "1.35%"
    "E17" = "-1.43%"
    "D18" = "3.167"
    "E18" = "-0.64%"
    "D19" = "2.176"
    "E19" = "-1.26%"
    "E20" = "2.41%"
    "E21" = "1.57%"
    "D22" = "4.094"
    "E22" = "7.31%"
    "D23" = "0.04254"
    "E23" = "-0.24%"
    "E24" = "-0.15%"
    "E25" = "-8.89%"
    "D26" = "0.0001201"
    "E26" = "0.09%"
    "E27" = "-18.89%"
    "D40" = "0.03883"
    "E40" = "1.37%"
    "E41" = "1.16%"
    "D42" = "0.004128"
    "E42" = "-34.76%"
    "E43" = "6.14%"
    "E44" = "0.42%"
    "D45" = "0.00005390"
    "E45" = "-0.21%"
    "E46" = "0.06%"
    "E47" = "7.00%"
    "D48" = "0.1321"
    "E48" = "-47.49%"
    "E49" = "0.06%"
    "E50" = "0.06%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
